$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hide columns C, D and I (their widths stay the same, only visibility changes)
$ws.Columns.Item(3).Hidden = $true
$ws.Columns.Item(4).Hidden = $true
$ws.Columns.Item(9).Hidden = $true

# Select column I and scroll the view back to the top-left corner (A1)
$ws.Range("I:I").Select() | Out-Null
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
